$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-24"

# Update the label for the August row
$ws.Range("A9").Value = "August (through 08-24)"

# Update August (row 9) figures
$ws.Range("B9").Value = 26
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = 68
$ws.Range("E9").Value = 43
$ws.Range("F9").Value = 35
$ws.Range("G9").Value = 138
$ws.Range("H9").Value = 126
$ws.Range("I9").Value = 132

# Update Total (row 10) figures
$ws.Range("B10").Value = 188
$ws.Range("C10").Value = 357
$ws.Range("D10").Value = 533
$ws.Range("E10").Value = 468
$ws.Range("F10").Value = 339
$ws.Range("G10").Value = 759
$ws.Range("H10").Value = 1036
$ws.Range("I10").Value = 1103
